# Update header text (volume/issue number and date range)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Volume 30   Number  13"
$ws.Range("C9").Value = "Report Covering the Week  3/27/2023  Through  4/2/2023"

# Update weekly crime statistics table (rows 15-27, 30)
# Row 15
$ws.Range("C14").Copy($ws.Range("C15"))
$ws.Range("E15").Value = -100
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = -75
$ws.Range("J15").Value = 7
$ws.Range("K15").Value = -28.571428571428
$ws.Range("N15").Value = -50

# Row 16
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 33.333333333333
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = 25
$ws.Range("I16").Value = 35
$ws.Range("J16").Value = 48
$ws.Range("K16").Value = -27.083333333333
$ws.Range("L16").Value = 84.210526315789
$ws.Range("M16").Value = -49.275362318840
$ws.Range("N16").Value = -86

# Row 17
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 33.333333333333
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 8
$ws.Range("H17").Value = 112.5
$ws.Range("I17").Value = 61
$ws.Range("J17").Value = 38
$ws.Range("K17").Value = 60.526315789473
$ws.Range("L17").Value = 74.285714285714
$ws.Range("M17").Value = 190.47619047619
$ws.Range("N17").Value = 1.666666666666

# Row 18
$ws.Range("C18").Value = 7
$ws.Range("E18").Value = 250
$ws.Range("F18").Value = 22
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = 57.142857142857
$ws.Range("I18").Value = 68
$ws.Range("J18").Value = 42
$ws.Range("K18").Value = 61.904761904761
$ws.Range("L18").Value = 61.904761904761
$ws.Range("M18").Value = 7.936507936507
$ws.Range("N18").Value = -82.871536523929

# Row 19
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -33.333333333333
$ws.Range("F19").Value = 37
$ws.Range("G19").Value = 39
$ws.Range("H19").Value = -5.128205128205
$ws.Range("I19").Value = 159
$ws.Range("J19").Value = 167
$ws.Range("K19").Value = -4.790419161676
$ws.Range("L19").Value = 123.943661971831
$ws.Range("M19").Value = 32.5
$ws.Range("N19").Value = 20.454545454545

# Row 20
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 20
$ws.Range("H20").Value = 17.647058823529
$ws.Range("I20").Value = 67
$ws.Range("J20").Value = 61
$ws.Range("K20").Value = 9.836065573770
$ws.Range("L20").Value = 252.631578947368
$ws.Range("M20").Value = 36.734693877551
$ws.Range("N20").Value = -94.517184942716

# Row 21
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = 25
$ws.Range("F21").Value = 112
$ws.Range("G21").Value = 94
$ws.Range("H21").Value = 19.148936170212
$ws.Range("I21").Value = 395
$ws.Range("J21").Value = 363
$ws.Range("K21").Value = 8.815426997245
$ws.Range("L21").Value = 110.106382978723
$ws.Range("M21").Value = 21.538461538461
$ws.Range("N21").Value = -80.945489628557

# Row 22
$ws.Range("D15").Copy($ws.Range("C22"))
$ws.Range("C22").Value = 1
$ws.Range("C14").Copy($ws.Range("D22"))
$ws.Range("E14").Copy($ws.Range("E22"))
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 33.333333333333
$ws.Range("I22").Value = 10
$ws.Range("K22").Value = 11.111111111111
$ws.Range("L22").Value = 25
$ws.Range("M22").Value = 66.666666666666

# Row 23
$ws.Range("C14").Copy($ws.Range("C23"))
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 5
$ws.Range("H23").Value = 66.666666666666
$ws.Range("J23").Value = 15
$ws.Range("K23").Value = 40

# Row 24
$ws.Range("C24").Value = 32
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = 52.380952380952
$ws.Range("F24").Value = 108
$ws.Range("G24").Value = 77
$ws.Range("H24").Value = 40.259740259740
$ws.Range("I24").Value = 384
$ws.Range("J24").Value = 292
$ws.Range("K24").Value = 31.506849315068
$ws.Range("L24").Value = 94.923857868020
$ws.Range("M24").Value = 108.695652173913

# Row 25
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = -12.5
$ws.Range("G25").Value = 29
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 124
$ws.Range("J25").Value = 101
$ws.Range("K25").Value = 22.772277227722
$ws.Range("L25").Value = 37.777777777777
$ws.Range("M25").Value = 29.166666666666

# Row 26
$ws.Range("C14").Copy($ws.Range("C26"))
$ws.Range("E26").Value = -100
$ws.Range("G26").Value = 5
$ws.Range("H26").Value = -80
$ws.Range("J26").Value = 9
$ws.Range("K26").Value = -22.222222222222

# Row 27
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 12
$ws.Range("J27").Value = 17
$ws.Range("K27").Value = -29.411764705882
$ws.Range("L27").Value = 71.428571428571

# Row 30
$ws.Range("D15").Copy($ws.Range("F30"))
$ws.Range("F30").Value = 3
$ws.Range("I30").Value = 5
$ws.Range("K30").Value = 400

